$d = $word.ActiveDocument

$replacements = @(
    @{old = "306×5=1530"; new = "655×8=5240"},
    @{old = "383×4=1532"; new = "415×3=1245"},
    @{old = "676×7=4732"; new = "713×9=6417"},
    @{old = "810×9=7290"; new = "181×8=1448"},
    @{old = "880×9=7920"; new = "441×7=3087"},
    @{old = "607×3=1821"; new = "541×9=4869"},
    @{old = "356×5=1780"; new = "349×6=2094"},
    @{old = "522×7=3654"; new = "550×4=2200"},
    @{old = "500×6=3000"; new = "758×3=2274"},
    @{old = "667×8=5336"; new = "940×4=3760"},
    @{old = "394×2=788";  new = "889×2=1778"},
    @{old = "644×5=3220"; new = "718×6=4308"},
    @{old = "189×8=1512"; new = "594×9=5346"},
    @{old = "920×3=2760"; new = "226×6=1356"},
    @{old = "129×3=387";  new = "974×7=6818"},
    @{old = "780×4=3120"; new = "125×4=500"},
    @{old = "681×4=2724"; new = "367×9=3303"},
    @{old = "894×8=7152"; new = "536×8=4288"},
    @{old = "686×3=2058"; new = "611×2=1222"},
    @{old = "509×3=1527"; new = "531×3=1593"},
    @{old = "879×9=7911"; new = "293×5=1465"},
    @{old = "457×9=4113"; new = "377×7=2639"},
    @{old = "196×2=392";  new = "665×9=5985"},
    @{old = "220×3=660";  new = "124×3=372"},
    @{old = "723×7=5061"; new = "655×8=5240"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "Done replacing $($replacements.Count) cells"
